# The "Recorded By" column (G) on the active sheet contains comma-separated
# lists of recorder names/emails (e.g. "System, dnasr281@gmail.com"). This
# edit reverses the order of the items in that comma-separated list for the
# rows below (the rows whose list order had not already been flipped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToReverse = @(
    2,3,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,26,28,29,31,32,33,34,36,37,38,39,40,
    41,43,44,45,46,47,48,50,52,54,55,57,58,59,60,62,63,64,65,66,67,69,70,71,72,73,74,76,
    78,80,81,82,83,84,85,86,87,90,92,93,94,96,99,101,106,107,108,109,110,111,112,113,116,
    118,119,120,122,125,127,132,133,134,135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($r in $rowsToReverse) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    $parts = $current -split ", "
    if ($parts.Count -gt 1) {
        $reversedParts = $parts[($parts.Count - 1)..0]
        $cell.Value = ($reversedParts -join ", ")
    }
}
